$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 353.42856
$ws.Range("I28").Value = 353.42856
$ws.Range("K28").Value = 353.42856
$ws.Range("M28").Value = 131.57144

$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 2000
$ws.Range("K32").Value = 2000
$ws.Range("M32").Value = -1674

$ws.Range("H55").Value = 715
$ws.Range("I55").Value = 423.5
$ws.Range("J55").Value = 1298
$ws.Range("K55").Value = 423.5
$ws.Range("L55").Value = 1298
$ws.Range("M55").Value = -209.5
$ws.Range("N55").Value = -1726

$ws.Range("H80").Value = 2717.3635
$ws.Range("I80").Value = 1249
$ws.Range("J80").Value = 3556.4285
$ws.Range("K80").Value = 3747
$ws.Range("L80").Value = 10669.2855
$ws.Range("M80").Value = -2749
$ws.Range("N80").Value = -12665.2855

$ws.Range("H83").Value = 2717.3635
$ws.Range("I83").Value = 1249
$ws.Range("J83").Value = 3556.4285
$ws.Range("K83").Value = 11241
$ws.Range("L83").Value = 32007.8565
$ws.Range("M83").Value = -6249
$ws.Range("N83").Value = -41991.8565

$ws.Range("H88").Value = 2531.125
$ws.Range("J88").Value = 2642.3333
$ws.Range("L88").Value = 2642.3333
$ws.Range("N88").Value = -3454.3333

$ws.Range("H91").Value = 2531.125
$ws.Range("J91").Value = 2642.3333
$ws.Range("L91").Value = 2642.3333
$ws.Range("N91").Value = -5450.3333

$ws.Range("H138").Value = 6479.967
$ws.Range("J138").Value = 6681.4443
$ws.Range("L138").Value = 20044.3329
$ws.Range("N138").Value = -30324.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9709.171
$ws.Range("I32").Value = 7655.795
$ws.Range("K32").Value = 7655.795
$ws.Range("M32").Value = -7368.795

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5503.1113
$ws.Range("I20").Value = 1298.6666
$ws.Range("K20").Value = 1298.6666
$ws.Range("M20").Value = -1051.6666

$ws.Range("H97").Value = 29725
$ws.Range("I97").Value = 29700
$ws.Range("J97").Value = 29750
$ws.Range("K97").Value = 29700
$ws.Range("L97").Value = 29750
$ws.Range("M97").Value = -28709
$ws.Range("N97").Value = -31732

$ws.Range("H100").Value = 80000
$ws.Range("J100").Value = 80000
$ws.Range("L100").Value = 80000
$ws.Range("N100").Value = -82164

$ws.Range("H103").Value = 78666.336
$ws.Range("J103").Value = 78666.336
$ws.Range("L103").Value = 78666.336
$ws.Range("N103").Value = -81010.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 191.9
$ws.Range("I7").Value = 134.4
$ws.Range("J7").Value = 249.4
$ws.Range("K7").Value = 134.4
$ws.Range("L7").Value = 249.4
$ws.Range("M7").Value = -21.40000000000001
$ws.Range("N7").Value = -475.4

$ws.Range("H58").Value = 4703.5713
$ws.Range("I58").Value = 2898.3333
$ws.Range("K58").Value = 2898.3333
$ws.Range("M58").Value = -2695.3333

$ws.Range("H86").Value = 7752.2666
$ws.Range("I86").Value = 6998.3335
$ws.Range("K86").Value = 6998.3335
$ws.Range("M86").Value = -5875.3335

$ws.Range("H89").Value = 7752.2666
$ws.Range("I89").Value = 6998.3335
$ws.Range("K89").Value = 34991.6675
$ws.Range("M89").Value = -29375.6675

$ws.Range("H99").Value = 3500.8823
$ws.Range("I99").Value = 3543.5715
$ws.Range("J99").Value = 3301.6667
$ws.Range("K99").Value = 3543.5715
$ws.Range("L99").Value = 3301.6667
$ws.Range("M99").Value = -2045.5715
$ws.Range("N99").Value = -6297.6667

$ws.Range("H126").Value = 3500.8823
$ws.Range("I126").Value = 3543.5715
$ws.Range("J126").Value = 3301.6667
$ws.Range("K126").Value = 10630.7145
$ws.Range("L126").Value = 9905.000100000001
$ws.Range("M126").Value = -8160.7145
$ws.Range("N126").Value = -14845.0001

$ws.Range("H132").Value = 3487.3333
$ws.Range("I132").Value = 2938.4
$ws.Range("K132").Value = 8815.200000000001
$ws.Range("M132").Value = -6285.200000000001

$ws.Range("H136").Value = 4703.5713
$ws.Range("I136").Value = 2898.3333
$ws.Range("K136").Value = 8694.999899999999
$ws.Range("M136").Value = -6144.999899999999

$ws.Range("H141").Value = 598993.5
$ws.Range("J141").Value = 698658
$ws.Range("L141").Value = 698658
$ws.Range("N141").Value = -709018

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 572.36365
$ws.Range("I12").Value = 1003.6667
$ws.Range("J12").Value = 410.625
$ws.Range("K12").Value = 3011.0001
$ws.Range("L12").Value = 1231.875
$ws.Range("M12").Value = -2838.0001
$ws.Range("N12").Value = -1577.875

$ws.Range("H131").Value = 4058.8235
$ws.Range("J131").Value = 4133.3335
$ws.Range("L131").Value = 12400.0005
$ws.Range("N131").Value = -22480.0005

$ws.Range("H140").Value = 3480
$ws.Range("I140").Value = 3480
$ws.Range("K140").Value = 10440
$ws.Range("M140").Value = -5260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5500
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 1000
$ws.Range("M70").Value = -730

$ws.Range("H73").Value = 5500
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 1000
$ws.Range("M73").Value = -64

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9364.444
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10590

$ws.Range("H27").Value = 9364.444
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10214

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H46").Value = 2112
$ws.Range("J46").Value = 2112
$ws.Range("L46").Value = 2112
$ws.Range("N46").Value = -2488

$ws.Range("H55").Value = 907.06665
$ws.Range("I55").Value = 557.875
$ws.Range("J55").Value = 1306.1428
$ws.Range("K55").Value = 557.875
$ws.Range("L55").Value = 1306.1428
$ws.Range("M55").Value = -384.875
$ws.Range("N55").Value = -1652.1428

$ws.Range("H82").Value = 2094.4167
$ws.Range("J82").Value = 2110.5
$ws.Range("L82").Value = 2110.5
$ws.Range("N82").Value = -2832.5

$ws.Range("H85").Value = 2094.4167
$ws.Range("J85").Value = 2110.5
$ws.Range("L85").Value = 2110.5
$ws.Range("N85").Value = -4606.5

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 3511.2
$ws.Range("I132").Value = 2889.25
$ws.Range("K132").Value = 8667.75
$ws.Range("M132").Value = -6137.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7270.7144
$ws.Range("I62").Value = 3367.3333
$ws.Range("K62").Value = 3367.3333
$ws.Range("M62").Value = -2743.3333

$ws.Range("H65").Value = 7270.7144
$ws.Range("I65").Value = 3367.3333
$ws.Range("K65").Value = 16836.6665
$ws.Range("M65").Value = -13716.6665

$ws.Range("H81").Value = 3206.6365
$ws.Range("I81").Value = 2627.3
$ws.Range("K81").Value = 5254.6
$ws.Range("M81").Value = -4193.6

$ws.Range("H84").Value = 3206.6365
$ws.Range("I84").Value = 2627.3
$ws.Range("K84").Value = 26273
$ws.Range("M84").Value = -20969

$ws.Range("H132").Value = 2875.4167
$ws.Range("I132").Value = 1542.7368
$ws.Range("K132").Value = 4628.2104
$ws.Range("M132").Value = -2098.2104
